# Trade #29 closed at 2026-02-17 20:54:51 - unknown UNKNOWN +0.000%
#
# This script updates the "live_trading_results" workbook:
#  - Refreshes aggregate metrics on the Summary and Strategy Status sheets
#  - Marks the previously OPEN MarketMaking trade (trade #57) as CLOSED
#    on both the "All Trades" and "MarketMaking" sheets
#  - Appends a brand-new OPEN MarketMaking trade (trade #90) to both
#    the "All Trades" and "MarketMaking" sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1400.6    # Current Capital
$wsSummary.Range("B4").Value = 0.39      # Total P&L $
$wsSummary.Range("B6").Value = 57        # Total Trades
$wsSummary.Range("B7").Value = 29        # Winning Trades
$wsSummary.Range("B9").Value = 50.88     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet (row 5 = MarketMaking)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 100.6   # Capital
$wsStatus.Range("D5").Value = 24      # Trades
$wsStatus.Range("E5").Value = 0.28    # P&L $
$wsStatus.Range("F5").Value = 0.6     # P&L %
$wsStatus.Range("G5").Value = 62.5    # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Trade #57 (row 58) transitions from OPEN -> CLOSED
$wsAll.Cells.Item(58, 7).Value = 0.98        # G: Exit Price
$wsAll.Cells.Item(58, 8).Value = "CLOSED"    # H: Status
$wsAll.Cells.Item(58, 9).Value = 1.0309      # I: P&L %
$wsAll.Cells.Item(58, 10).Value = 0.01       # J: P&L $
$wsAll.Cells.Item(58, 11).Value = 100.6      # K: Capital After
$wsAll.Cells.Item(58, 12).Value = "early_exit"  # L: Exit Reason
$wsAll.Cells.Item(58, 13).Value = 0.11       # M: Duration (min)

# New trade #90 (row 91), still OPEN
$wsAll.Cells.Item(91, 1).Value = 90           # A: Trade #

$wsAll.Cells.Item(91, 2).NumberFormat = "@"
$wsAll.Cells.Item(91, 2).Value = "2026-02-17" # B: Date
$wsAll.Cells.Item(91, 2).Style = "Normal"

$wsAll.Cells.Item(91, 3).NumberFormat = "@"
$wsAll.Cells.Item(91, 3).Value = "20:54:44"   # C: Time
$wsAll.Cells.Item(91, 3).Style = "Normal"

$wsAll.Cells.Item(91, 4).Value = "MarketMaking"  # D: Strategy
$wsAll.Cells.Item(91, 5).Value = "DOWN"          # E: Side
$wsAll.Cells.Item(91, 6).Value = 0.97             # F: Entry Price
# G: Exit Price left blank (still open)
$wsAll.Cells.Item(91, 8).Value = "OPEN"           # H: Status
$wsAll.Cells.Item(91, 9).Value = 0                # I: P&L %
$wsAll.Cells.Item(91, 10).Value = 0               # J: P&L $
$wsAll.Cells.Item(91, 11).Value = 100.5855022889912  # K: Capital After
# L: Exit Reason left blank
$wsAll.Cells.Item(91, 13).Value = 0               # M: Duration (min)
$wsAll.Cells.Item(91, 14).Value = 0               # N: Entry Slippage (bps)
$wsAll.Cells.Item(91, 15).Value = 0               # O: Exit Slippage (bps)
$wsAll.Cells.Item(91, 16).Value = 0.6             # P: Confidence
$wsAll.Cells.Item(91, 17).Value = "Normal spread capture: 19600 bps"  # Q: Entry Reason

# ---------------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

# Trade #57 (row 25) transitions from OPEN -> CLOSED
$wsMM.Cells.Item(25, 7).Value = 0.98        # G: Exit Price
$wsMM.Cells.Item(25, 8).Value = "CLOSED"    # H: Status
$wsMM.Cells.Item(25, 9).Value = 1.0309      # I: P&L %
$wsMM.Cells.Item(25, 10).Value = 0.01       # J: P&L $
$wsMM.Cells.Item(25, 11).Value = 100.6      # K: Capital After
$wsMM.Cells.Item(25, 16).Value = "early_exit"  # P: Exit Reason
$wsMM.Cells.Item(25, 17).Value = 0.11       # Q: Duration (min)

# New trade #90 (row 58), still OPEN
$wsMM.Cells.Item(58, 1).Value = 90            # A: Trade #

$wsMM.Cells.Item(58, 2).NumberFormat = "@"
$wsMM.Cells.Item(58, 2).Value = "2026-02-17"  # B: Date
$wsMM.Cells.Item(58, 2).Style = "Normal"

$wsMM.Cells.Item(58, 3).NumberFormat = "@"
$wsMM.Cells.Item(58, 3).Value = "20:54:44"    # C: Time
$wsMM.Cells.Item(58, 3).Style = "Normal"

$wsMM.Cells.Item(58, 4).Value = "MarketMaking"  # D: Strategy
$wsMM.Cells.Item(58, 5).Value = "DOWN"          # E: Side
$wsMM.Cells.Item(58, 6).Value = 0.97             # F: Entry Price
# G: Exit Price left blank (still open)
$wsMM.Cells.Item(58, 8).Value = "OPEN"           # H: Status
$wsMM.Cells.Item(58, 9).Value = 0                # I: P&L %
$wsMM.Cells.Item(58, 10).Value = 0               # J: P&L $
$wsMM.Cells.Item(58, 11).Value = 100.5855022889912  # K: Capital After
$wsMM.Cells.Item(58, 12).Value = 0               # L: Entry Slippage (bps)
$wsMM.Cells.Item(58, 13).Value = 0               # M: Exit Slippage (bps)
$wsMM.Cells.Item(58, 14).Value = 0.6             # N: Confidence
$wsMM.Cells.Item(58, 15).Value = "Normal spread capture: 19600 bps"  # O: Entry Reason
# P: Exit Reason left blank
$wsMM.Cells.Item(58, 17).Value = 0               # Q: Duration (min)
